$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19332.334
$ws.Range("I21").Value = 17999
$ws.Range("J21").Value = 19999
$ws.Range("K21").Value = 17999
$ws.Range("L21").Value = 19999
$ws.Range("M21").Value = -17531
$ws.Range("N21").Value = -20935
$ws.Range("H23").Value = 19332.334
$ws.Range("I23").Value = 17999
$ws.Range("J23").Value = 19999
$ws.Range("K23").Value = 17999
$ws.Range("L23").Value = 19999
$ws.Range("M23").Value = -17765
$ws.Range("N23").Value = -20467
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -19
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 424.3
$ws.Range("I38").Value = 223.28572
$ws.Range("J38").Value = 893.3333
$ws.Range("K38").Value = 669.85716
$ws.Range("L38").Value = 2679.9999
$ws.Range("M38").Value = -297.85716
$ws.Range("N38").Value = -3423.9999
$ws.Range("H58").Value = 1248.5
$ws.Range("I58").Value = 882.7692
$ws.Range("J58").Value = 2833.3333
$ws.Range("K58").Value = 2648.3076
$ws.Range("L58").Value = 8499.999899999999
$ws.Range("M58").Value = -2498.3076
$ws.Range("N58").Value = -8799.999899999999
$ws.Range("H87").Value = 24885.428
$ws.Range("J87").Value = 24885.428
$ws.Range("L87").Value = 24885.428
$ws.Range("N87").Value = -27381.428
$ws.Range("H90").Value = 24885.428
$ws.Range("J90").Value = 24885.428
$ws.Range("L90").Value = 74656.284
$ws.Range("N90").Value = -87136.284
$ws.Range("H92").Value = 11906048
$ws.Range("I92").Value = 12821628
$ws.Range("J92").Value = 3500
$ws.Range("K92").Value = 12821628
$ws.Range("L92").Value = 3500
$ws.Range("M92").Value = -12820380
$ws.Range("N92").Value = -5996
$ws.Range("H132").Value = 4221.6943
$ws.Range("I132").Value = 4027.1035
$ws.Range("J132").Value = 5027.857
$ws.Range("K132").Value = 12081.3105
$ws.Range("L132").Value = 15083.571
$ws.Range("M132").Value = -9551.3105
$ws.Range("N132").Value = -20143.571
$ws.Range("H137").Value = 3489.7778
$ws.Range("I137").Value = 2451.5356
$ws.Range("K137").Value = 7354.6068
$ws.Range("M137").Value = -4804.6068
$ws.Range("H138").Value = 1926.26
$ws.Range("I138").Value = 611.8125
$ws.Range("J138").Value = 2176.6309
$ws.Range("K138").Value = 1835.4375
$ws.Range("L138").Value = 6529.8927
$ws.Range("M138").Value = 3304.5625
$ws.Range("N138").Value = -16809.8927
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 13773.667
$ws.Range("J23").Value = 9757.071
$ws.Range("L23").Value = 9757.071
$ws.Range("N23").Value = -10275.071
$ws.Range("H32").Value = 4089.45
$ws.Range("I32").Value = 3524.7144
$ws.Range("J32").Value = 9799.556
$ws.Range("K32").Value = 3524.7144
$ws.Range("L32").Value = 9799.556
$ws.Range("M32").Value = -3237.7144
$ws.Range("N32").Value = -10373.556
$ws.Range("H37").Value = 10999.5
$ws.Range("H44").Value = 19312.572
$ws.Range("J44").Value = 19312.572
$ws.Range("L44").Value = 19312.572
$ws.Range("N44").Value = -20288.572
$ws.Range("H63").Value = 4203.4165
$ws.Range("I63").Value = 2938.111
$ws.Range("J63").Value = 7999.3335
$ws.Range("K63").Value = 2938.111
$ws.Range("L63").Value = 7999.3335
$ws.Range("M63").Value = -2252.111
$ws.Range("N63").Value = -9371.333500000001
$ws.Range("H66").Value = 4203.4165
$ws.Range("I66").Value = 2938.111
$ws.Range("J66").Value = 7999.3335
$ws.Range("K66").Value = 14690.555
$ws.Range("L66").Value = 39996.6675
$ws.Range("M66").Value = -11258.555
$ws.Range("N66").Value = -46860.6675
$ws.Range("H110").Value = 1174.8
$ws.Range("I110").Value = 1199.7778
$ws.Range("K110").Value = 1199.7778
$ws.Range("M110").Value = 845.2221999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15120.333
$ws.Range("I82").Value = 6135.6665
$ws.Range("J82").Value = 21110.111
$ws.Range("K82").Value = 6135.6665
$ws.Range("L82").Value = 21110.111
$ws.Range("M82").Value = -5752.6665
$ws.Range("N82").Value = -21876.111
$ws.Range("H85").Value = 15120.333
$ws.Range("I85").Value = 6135.6665
$ws.Range("J85").Value = 21110.111
$ws.Range("K85").Value = 6135.6665
$ws.Range("L85").Value = 21110.111
$ws.Range("M85").Value = -4809.6665
$ws.Range("N85").Value = -23762.111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2295
$ws.Range("I16").Value = 2295
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2295
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2008
$ws.Range("N16").ClearContents()
$ws.Range("H113").Value = 2295
$ws.Range("I113").Value = 2295
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2295
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -125
$ws.Range("N113").ClearContents()
$ws.Range("H134").Value = 2163.7646
$ws.Range("I134").Value = 2135.818
$ws.Range("J134").Value = 2215
$ws.Range("K134").Value = 6407.454000000001
$ws.Range("L134").Value = 6645
$ws.Range("M134").Value = -3872.454000000001
$ws.Range("N134").Value = -11715
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 508.2353
$ws.Range("I5").Value = 477.5
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1432.5
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1320.5
$ws.Range("N5").Value = -3224
$ws.Range("H23").Value = 41666776
$ws.Range("I23").Value = 122.666664
$ws.Range("J23").Value = 55555660
$ws.Range("K23").Value = 367.999992
$ws.Range("L23").Value = 166666980
$ws.Range("M23").Value = -132.999992
$ws.Range("N23").Value = -166667450
$ws.Range("H55").Value = 1745.4546
$ws.Range("J55").Value = 1745.4546
$ws.Range("L55").Value = 5236.3638
$ws.Range("N55").Value = -5590.3638
$ws.Range("H113").Value = 3072.5
$ws.Range("I113").Value = 880
$ws.Range("J113").Value = 3803.3333
$ws.Range("K113").Value = 2640
$ws.Range("L113").Value = 11409.9999
$ws.Range("M113").Value = -470
$ws.Range("N113").Value = -15749.9999
$ws.Range("H131").Value = 985.07275
$ws.Range("I131").Value = 395
$ws.Range("J131").Value = 1031.3529
$ws.Range("K131").Value = 1185
$ws.Range("L131").Value = 3094.0587
$ws.Range("M131").Value = 3855
$ws.Range("N131").Value = -13174.0587
$ws.Range("H132").Value = 2091.5173
$ws.Range("J132").Value = 1965.3914
$ws.Range("L132").Value = 17688.5226
$ws.Range("N132").Value = -22748.5226
$ws.Range("H135").Value = 508.2353
$ws.Range("I135").Value = 477.5
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 4297.5
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -1762.5
$ws.Range("N135").Value = -14070
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 5365
$ws.Range("I46").Value = 15515.5
$ws.Range("J46").Value = 4349.95
$ws.Range("K46").Value = 15515.5
$ws.Range("L46").Value = 4349.95
$ws.Range("M46").Value = -15359.5
$ws.Range("N46").Value = -4661.95
$ws.Range("H70").Value = 5586.1797
$ws.Range("I70").Value = 5597.6924
$ws.Range("J70").Value = 5563.154
$ws.Range("K70").Value = 5597.6924
$ws.Range("L70").Value = 5563.154
$ws.Range("M70").Value = -5327.6924
$ws.Range("N70").Value = -6103.154
$ws.Range("H73").Value = 5586.1797
$ws.Range("I73").Value = 5597.6924
$ws.Range("J73").Value = 5563.154
$ws.Range("K73").Value = 5597.6924
$ws.Range("L73").Value = 5563.154
$ws.Range("M73").Value = -4661.6924
$ws.Range("N73").Value = -7435.154
$ws.Range("H80").Value = 56602892
$ws.Range("I80").Value = 72717576
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 72717576
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -72716578
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 56602892
$ws.Range("I83").Value = 72717576
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 363587880
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -363582888
$ws.Range("N83").Value = -1017484
$ws.Range("H132").Value = 3692.1428
$ws.Range("I132").Value = 3635.3333
$ws.Range("J132").Value = 3734.75
$ws.Range("K132").Value = 10905.9999
$ws.Range("L132").Value = 11204.25
$ws.Range("M132").Value = -8375.999899999999
$ws.Range("N132").Value = -16264.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8282.346
$ws.Range("I22").Value = 764
$ws.Range("J22").Value = 10537.85
$ws.Range("K22").Value = 764
$ws.Range("L22").Value = 10537.85
$ws.Range("M22").Value = -469
$ws.Range("N22").Value = -11127.85
$ws.Range("H24").Value = 70007
$ws.Range("J24").Value = 70007
$ws.Range("L24").Value = 70007
$ws.Range("N24").Value = -70693
$ws.Range("H27").Value = 8282.346
$ws.Range("I27").Value = 764
$ws.Range("J27").Value = 10537.85
$ws.Range("K27").Value = 764
$ws.Range("L27").Value = 10537.85
$ws.Range("M27").Value = -657
$ws.Range("N27").Value = -10751.85
$ws.Range("H36").Value = 99715
$ws.Range("J36").Value = 99715
$ws.Range("L36").Value = 99715
$ws.Range("N36").Value = -100839
$ws.Range("H61").Value = 4167.8
$ws.Range("I61").Value = 3596.4167
$ws.Range("J61").Value = 6453.3335
$ws.Range("K61").Value = 3596.4167
$ws.Range("L61").Value = 6453.3335
$ws.Range("M61").Value = -3394.4167
$ws.Range("N61").Value = -6857.3335
$ws.Range("H113").Value = 4167.8
$ws.Range("I113").Value = 3596.4167
$ws.Range("J113").Value = 6453.3335
$ws.Range("K113").Value = 3596.4167
$ws.Range("L113").Value = 6453.3335
$ws.Range("M113").Value = -1426.4167
$ws.Range("N113").Value = -10793.3335
$ws.Range("H136").Value = 6946530.5
$ws.Range("I136").Value = 1667
$ws.Range("J136").Value = 13891394
$ws.Range("K136").Value = 5001
$ws.Range("L136").Value = 41674182
$ws.Range("M136").Value = -2451
$ws.Range("N136").Value = -41679282
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9999
$ws.Range("J54").Value = 9999
$ws.Range("L54").Value = 9999
$ws.Range("N54").Value = -11039
$ws.Range("H81").Value = 5546.9
$ws.Range("I81").Value = 6314.2
$ws.Range("J81").Value = 4779.6
$ws.Range("K81").Value = 12628.4
$ws.Range("L81").Value = 9559.200000000001
$ws.Range("M81").Value = -11567.4
$ws.Range("N81").Value = -11681.2
$ws.Range("H84").Value = 5546.9
$ws.Range("I84").Value = 6314.2
$ws.Range("J84").Value = 4779.6
$ws.Range("K84").Value = 63142
$ws.Range("L84").Value = 47796
$ws.Range("M84").Value = -57838
$ws.Range("N84").Value = -58404
$ws.Range("H132").Value = 5053228.5
$ws.Range("I132").Value = 4145.4
$ws.Range("J132").Value = 7248482
$ws.Range("K132").Value = 12436.2
$ws.Range("L132").Value = 21745446
$ws.Range("M132").Value = -9906.199999999999
$ws.Range("N132").Value = -21750506
